$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.405.44"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.25%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.710.14"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.17%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "224.57"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5327"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.97%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2673"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.97%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06629"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.53%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.96"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.97%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07622"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.91%  "
$ws.Range("E12").Value = "  -2.86%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.709.43"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.39%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.944.95"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.21%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5779"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.78%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅8191"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.26%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.85"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.51%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "27.383.72"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.30%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "218.93"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.93%  "
$ws.Range("E20").Value = "  +0.01%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.657"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.95%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.46"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.62%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.954"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.94%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.004"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "142.10"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.57%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.723"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1212"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.93%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.260"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.61%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "16.25"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.50%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05402"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.67%  "
$ws.Range("E31").Value = "  -1.41%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.500"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.19%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.426"
$ws.Range("D33").Style = "Normal"
$ws.Range("E34").Value = "  -1.67%  "
$ws.Range("E35").Value = "  +1.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9479"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.41%  "
$ws.Range("E37").Value = "  -0.98%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5874"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.71%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01637"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.93%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.860"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.76%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.047.18"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.13%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.004"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.10%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8426"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.31%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.98"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.51%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.852.83"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.09%  "
$ws.Range("E46").Value = "  +2.41%  "
$ws.Range("E47").Value = "  -2.58%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4517"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.95%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.004"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.47%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.067"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.71%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05225"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.95%  "
